$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 2961.4285
$ws.Range("I53").Value = 265
$ws.Range("J53").Value = 3410.8333
$ws.Range("K53").Value = 265
$ws.Range("L53").Value = 3410.8333
$ws.Range("M53").Value = 372
$ws.Range("N53").Value = -4684.8333
$ws.Range("H76").Value = 3971064.2
$ws.Range("J76").Value = 9261816
$ws.Range("L76").Value = 9261816
$ws.Range("N76").Value = -9262446
$ws.Range("H79").Value = 3971064.2
$ws.Range("J79").Value = 9261816
$ws.Range("L79").Value = 9261816
$ws.Range("N79").Value = -9264000
$ws.Range("H129").Value = 303940.5
$ws.Range("J129").Value = 417819.6
$ws.Range("L129").Value = 1253458.8
$ws.Range("N129").Value = -1263458.8
$ws.Range("H132").Value = 18520394
$ws.Range("I132").Value = 20410186
$ws.Range("J132").Value = 413.8
$ws.Range("K132").Value = 61230558
$ws.Range("L132").Value = 1241.4
$ws.Range("M132").Value = -61228028
$ws.Range("N132").Value = -6301.4
$ws.Range("H138").Value = 3938.1736
$ws.Range("I138").Value = 3382.1667
$ws.Range("J138").Value = 4063.275
$ws.Range("K138").Value = 10146.5001
$ws.Range("L138").Value = 12189.825
$ws.Range("M138").Value = -5006.500100000001
$ws.Range("N138").Value = -22469.825
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 16733889
$ws.Range("J61").Value = 4320
$ws.Range("L61").Value = 4320
$ws.Range("N61").Value = -4744
$ws.Range("H74").Value = 31253088
$ws.Range("I74").Value = 50002784
$ws.Range("K74").Value = 50002784
$ws.Range("M74").Value = -50001910
$ws.Range("H77").Value = 31253088
$ws.Range("I77").Value = 50002784
$ws.Range("K77").Value = 250013920
$ws.Range("M77").Value = -250009552
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H136").Value = 16733889
$ws.Range("J136").Value = 4320
$ws.Range("L136").Value = 12960
$ws.Range("N136").Value = -18060
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H113").Value = 3000
$ws.Range("I113").Value = 3000
$ws.Range("K113").Value = 3000
$ws.Range("M113").Value = -830
$ws.Range("H128").Value = 3039
$ws.Range("I128").Value = 3039
$ws.Range("K128").Value = 9117
$ws.Range("M128").Value = -6627
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7665.3335
$ws.Range("I31").Value = 4478
$ws.Range("J31").Value = 8764.414000000001
$ws.Range("K31").Value = 4478
$ws.Range("L31").Value = 8764.414000000001
$ws.Range("M31").Value = -4183
$ws.Range("N31").Value = -9354.414000000001
$ws.Range("H34").Value = 7665.3335
$ws.Range("I34").Value = 4478
$ws.Range("J34").Value = 8764.414000000001
$ws.Range("K34").Value = 4478
$ws.Range("L34").Value = 8764.414000000001
$ws.Range("M34").Value = -4276
$ws.Range("N34").Value = -9168.414000000001
$ws.Range("H76").Value = 25002250
$ws.Range("I76").Value = 25002250
$ws.Range("K76").Value = 25002250
$ws.Range("M76").Value = -25001935
$ws.Range("H79").Value = 25002250
$ws.Range("I79").Value = 25002250
$ws.Range("K79").Value = 25002250
$ws.Range("M79").Value = -25001158
$ws.Range("H86").Value = 8177.579
$ws.Range("I86").Value = 2040.1
$ws.Range("K86").Value = 2040.1
$ws.Range("M86").Value = -917.0999999999999
$ws.Range("H89").Value = 8177.579
$ws.Range("I89").Value = 2040.1
$ws.Range("K89").Value = 10200.5
$ws.Range("M89").Value = -4584.5
$ws.Range("H122").Value = 4884.4287
$ws.Range("I122").Value = 6375.25
$ws.Range("J122").Value = 2896.6667
$ws.Range("K122").Value = 19125.75
$ws.Range("L122").Value = 8690.000100000001
$ws.Range("M122").Value = -16675.75
$ws.Range("N122").Value = -13590.0001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2780
$ws.Range("I132").Value = 833.3333
$ws.Range("J132").Value = 5700
$ws.Range("K132").Value = 7499.9997
$ws.Range("L132").Value = 51300
$ws.Range("M132").Value = -4969.9997
$ws.Range("N132").Value = -56360
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2504101.2
$ws.Range("I70").Value = 4135.2666
$ws.Range("J70").Value = 6254050
$ws.Range("K70").Value = 4135.2666
$ws.Range("L70").Value = 6254050
$ws.Range("M70").Value = -3865.2666
$ws.Range("N70").Value = -6254590
$ws.Range("H73").Value = 2504101.2
$ws.Range("I73").Value = 4135.2666
$ws.Range("J73").Value = 6254050
$ws.Range("K73").Value = 4135.2666
$ws.Range("L73").Value = 6254050
$ws.Range("M73").Value = -3199.2666
$ws.Range("N73").Value = -6255922
$ws.Range("H132").Value = 4401099
$ws.Range("I132").Value = 7060702
$ws.Range("K132").Value = 21182106
$ws.Range("M132").Value = -21179576
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4167.6665
$ws.Range("I22").Value = 5250.5
$ws.Range("J22").Value = 2002
$ws.Range("K22").Value = 5250.5
$ws.Range("L22").Value = 2002
$ws.Range("M22").Value = -4955.5
$ws.Range("N22").Value = -2592
$ws.Range("H27").Value = 4167.6665
$ws.Range("I27").Value = 5250.5
$ws.Range("J27").Value = 2002
$ws.Range("K27").Value = 5250.5
$ws.Range("L27").Value = 2002
$ws.Range("M27").Value = -5143.5
$ws.Range("N27").Value = -2216
$ws.Range("H68").Value = 2455.5
$ws.Range("I68").Value = 1800
$ws.Range("J68").Value = 2848.8
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 2848.8
$ws.Range("M68").Value = -1051
$ws.Range("N68").Value = -4346.8
$ws.Range("H71").Value = 2455.5
$ws.Range("I71").Value = 1800
$ws.Range("J71").Value = 2848.8
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 14244
$ws.Range("M71").Value = -5256
$ws.Range("N71").Value = -21732
$ws.Range("H107").Value = 1560
$ws.Range("I107").Value = 1560
$ws.Range("K107").Value = 1560
$ws.Range("M107").Value = 360
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4372.5
$ws.Range("I62").Value = 4330
$ws.Range("J62").Value = 4500
$ws.Range("K62").Value = 4330
$ws.Range("L62").Value = 4500
$ws.Range("M62").Value = -3706
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 4372.5
$ws.Range("I65").Value = 4330
$ws.Range("J65").Value = 4500
$ws.Range("K65").Value = 21650
$ws.Range("L65").Value = 22500
$ws.Range("M65").Value = -18530
$ws.Range("N65").Value = -28740
